$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 813.55
$ws.Range("J17").Value = 813.55
$ws.Range("L17").Value = 2440.65
$ws.Range("N17").Value = -2776.65

$ws.Range("H53").Value = 403.8
$ws.Range("I53").Value = 151.42857
$ws.Range("K53").Value = 151.42857
$ws.Range("M53").Value = 485.57143

$ws.Range("H125").Value = 4206.7
$ws.Range("I125").Value = 4792
$ws.Range("K125").Value = 43128
$ws.Range("M125").Value = -40668

$ws.Range("H138").Value = 3333.8
$ws.Range("I138").Value = 3551
$ws.Range("J138").Value = 3231.5881
$ws.Range("K138").Value = 10653
$ws.Range("L138").Value = 9694.764299999999
$ws.Range("M138").Value = -5513
$ws.Range("N138").Value = -19974.7643

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2754.5134
$ws.Range("I63").Value = 4544.375
$ws.Range("K63").Value = 4544.375
$ws.Range("M63").Value = -3858.375

$ws.Range("H66").Value = 2754.5134
$ws.Range("I66").Value = 4544.375
$ws.Range("K66").Value = 22721.875
$ws.Range("M66").Value = -19289.875

$ws.Range("H102").Value = 1837.7778
$ws.Range("I102").Value = 1914.8667
$ws.Range("K102").Value = 1914.8667
$ws.Range("M102").Value = -292.8667

$ws.Range("H132").Value = 47013.523
$ws.Range("I132").Value = 47013.523
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 141040.569
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -138510.569
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 3030.5386
$ws.Range("I64").Value = 2414.2856
$ws.Range("J64").Value = 3749.5
$ws.Range("K64").Value = 2414.2856
$ws.Range("L64").Value = 3749.5
$ws.Range("M64").Value = -2189.2856
$ws.Range("N64").Value = -4199.5

$ws.Range("H67").Value = 3030.5386
$ws.Range("I67").Value = 2414.2856
$ws.Range("J67").Value = 3749.5
$ws.Range("K67").Value = 2414.2856
$ws.Range("L67").Value = 3749.5
$ws.Range("M67").Value = -1634.2856
$ws.Range("N67").Value = -5309.5

$ws.Range("H99").Value = 118973.11
$ws.Range("I99").Value = 502505
$ws.Range("K99").Value = 502505
$ws.Range("M99").Value = -501007

$ws.Range("H105").Value = 2945
$ws.Range("I105").Value = 3071.9
$ws.Range("K105").Value = 3071.9
$ws.Range("M105").Value = -1324.9

$ws.Range("H130").Value = 79998.5
$ws.Range("J130").Value = 79998.5
$ws.Range("L130").Value = 79998.5
$ws.Range("N130").Value = -90038.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 377
$ws.Range("I19").Value = 363.2
$ws.Range("K19").Value = 363.2
$ws.Range("M19").Value = -193.2

$ws.Range("H24").Value = 377
$ws.Range("I24").Value = 363.2
$ws.Range("K24").Value = 363.2
$ws.Range("M24").Value = -193.2

$ws.Range("H99").Value = 2482.9
$ws.Range("I99").Value = 1918.6
$ws.Range("J99").Value = 3047.2
$ws.Range("K99").Value = 1918.6
$ws.Range("L99").Value = 3047.2
$ws.Range("M99").Value = -420.5999999999999
$ws.Range("N99").Value = -6043.2

$ws.Range("H122").Value = 5772.5713
$ws.Range("I122").Value = 2012
$ws.Range("K122").Value = 6036
$ws.Range("M122").Value = -3586

$ws.Range("H126").Value = 2482.9
$ws.Range("I126").Value = 1918.6
$ws.Range("J126").Value = 3047.2
$ws.Range("K126").Value = 5755.799999999999
$ws.Range("L126").Value = 9141.599999999999
$ws.Range("M126").Value = -3285.799999999999
$ws.Range("N126").Value = -14081.6

$ws.Range("H132").Value = 2599.75
$ws.Range("I132").Value = 2599.75
$ws.Range("K132").Value = 7799.25
$ws.Range("M132").Value = -5269.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 8785.571
$ws.Range("J19").Value = 8785.571
$ws.Range("L19").Value = 26356.713
$ws.Range("N19").Value = -26704.713

$ws.Range("H114").Value = 897.65
$ws.Range("I114").Value = 932
$ws.Range("J114").Value = 846.125
$ws.Range("K114").Value = 2796
$ws.Range("L114").Value = 2538.375
$ws.Range("M114").Value = 458
$ws.Range("N114").Value = -9046.375

$ws.Range("H123").Value = 28837.166
$ws.Range("J123").Value = 13331.333
$ws.Range("L123").Value = 39993.999
$ws.Range("N123").Value = -44893.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1232.0625
$ws.Range("I97").Value = 947.3077
$ws.Range("J97").Value = 2466
$ws.Range("K97").Value = 947.3077
$ws.Range("L97").Value = 2466
$ws.Range("M97").Value = -451.3077
$ws.Range("N97").Value = -3458

$ws.Range("H102").Value = 4125.067
$ws.Range("J102").Value = 4426
$ws.Range("L102").Value = 4426
$ws.Range("N102").Value = -7670

$ws.Range("H122").Value = 3952.5
$ws.Range("I122").Value = 2940.2856
$ws.Range("J122").Value = 5369.6
$ws.Range("K122").Value = 8820.856800000001
$ws.Range("L122").Value = 16108.8
$ws.Range("M122").Value = -6370.856800000001
$ws.Range("N122").Value = -21008.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 16006
$ws.Range("I13").Value = 16006
$ws.Range("K13").Value = 16006
$ws.Range("M13").Value = -15866

$ws.Range("H55").Value = 669.2069
$ws.Range("I55").Value = 605.9474
$ws.Range("J55").Value = 789.4
$ws.Range("K55").Value = 605.9474
$ws.Range("L55").Value = 789.4
$ws.Range("M55").Value = -432.9474
$ws.Range("N55").Value = -1135.4

$ws.Range("H68").Value = 4165.864
$ws.Range("I68").Value = 2526.8
$ws.Range("K68").Value = 2526.8
$ws.Range("M68").Value = -1777.8

$ws.Range("H71").Value = 4165.864
$ws.Range("I71").Value = 2526.8
$ws.Range("K71").Value = 12634
$ws.Range("M71").Value = -8890

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H93").Value = 2629.3901
$ws.Range("I93").Value = 2355.2122
$ws.Range("J93").Value = 3760.375
$ws.Range("K93").Value = 2355.2122
$ws.Range("L93").Value = 3760.375
$ws.Range("M93").Value = -1107.2122
$ws.Range("N93").Value = -6256.375

$ws.Range("H132").Value = 52878.625
$ws.Range("I132").Value = 72527.47
$ws.Range("K132").Value = 217582.41
$ws.Range("M132").Value = -215052.41

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 250
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H19").Value = 300
$ws.Range("J19").Value = 300
$ws.Range("L19").Value = 300
$ws.Range("M19").Value = -648

$ws.Range("H100").Value = 1264.7273
$ws.Range("I100").Value = 1078.7142
$ws.Range("K100").Value = 2157.4284
$ws.Range("M100").Value = -1616.4284

$ws.Range("H107").Value = 2857.8
$ws.Range("I107").Value = 1430.5
$ws.Range("K107").Value = 4291.5
$ws.Range("M107").Value = -2371.5

$ws.Range("H125").Value = 67016
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

$ws.Range("H132").Value = 25889.62
$ws.Range("I132").Value = 25889.62
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 77668.86
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -75138.86
$ws.Range("N132").ClearContents()
